$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("40:41").Insert()

$ws.Range("A40").Value = "sound effect volume"
$ws.Range("B40").Value = "Sound Effect Volume"
$ws.Range("A41").Value = "music volume"
$ws.Range("B41").Value = "Music Volume"

for ($r = 38; $r -le 44; $r++) {
    $a = $ws.Cells.Item($r, 1).Text
    $b = $ws.Cells.Item($r, 2).Text
    Write-Output ("Row $r -> A='$a' B='$b'")
}
Write-Output ("Dimension: " + $ws.UsedRange.Address())
